# Replace the final "Socks In the Dark" paragraph (which contained the
# literal "2) " run and the two-goal blurb) with:
#   1. that same blurb text, now as an actual numbered-list item (numId 2)
#      instead of a manually-typed "2) " + ind-left indent, and
#   2. a brand-new paragraph giving the "possible solution" for goal (a),
#      also part of the numId 2 list, carrying the bookmark that used to
#      sit at the end of the original paragraph.

$d = $word.ActiveDocument

$target = $d.Paragraphs.Last
if ($target.Range.Text -notlike "2) So we have two goals*") {
    throw "Unexpected target paragraph text: $($target.Range.Text)"
}

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="2"/>
</w:numPr>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t xml:space="preserve">So we have two goals in this one. Have at least one matching pair and then have a matching pair of each color. </w:t>
</w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="2"/>
</w:numPr>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t xml:space="preserve">Possible </w:t>
</w:r>
<w:r>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t>solution</w:t>
</w:r>
<w:r>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t xml:space="preserve"> for goal </w:t>
</w:r>
<w:r>
<w:rPr>
<w:b/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t xml:space="preserve">a: </w:t>
</w:r>
<w:r>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t>Since there are three different colors and you only need one matching pair</w:t>
</w:r>
<w:r>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t>;</w:t>
</w:r>
<w:r>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:r>
<w:rPr>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
</w:rPr>
<w:t>the least amount of socks you will need to grab is 4. However this solution will not work for goal b.</w:t>
</w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
'@

$target.Range.InsertXML($newXml)
